$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 44881.50060109999
$ws.Range("A3").Value = 44881.50268443332
$ws.Range("A4").Value = 44881.50476776665
$ws.Range("A5").Value = 44881.50685109999
$ws.Range("A6").Value = 44881.50893443332
$ws.Range("A7").Value = 44881.51101776666
$ws.Range("A8").Value = 44881.51310109999
$ws.Range("A9").Value = 44881.51518443332
$ws.Range("A10").Value = 44881.51726776666
$ws.Range("A11").Value = 44881.51935109999
$ws.Range("A12").Value = 44881.52143443332
$ws.Range("A13").Value = 44881.52351776665
$ws.Range("A14").Value = 44881.52560109999
$ws.Range("A15").Value = 44881.52768443332
$ws.Range("A16").Value = 44881.52976776665
$ws.Range("A17").Value = 44881.53185109999
$ws.Range("A18").Value = 44881.53393443332
$ws.Range("A19").Value = 44881.53601776665
$ws.Range("A20").Value = 44881.53810109999
$ws.Range("A21").Value = 44881.54018443332
$ws.Range("A22").Value = 44881.54226776666
$ws.Range("A23").Value = 44881.54435109999
$ws.Range("A24").Value = 44881.54643443332
$ws.Range("A25").Value = 44881.54851776666
$ws.Range("A26").Value = 44881.55060109999
$ws.Range("A27").Value = 44881.55268443332
$ws.Range("A28").Value = 44881.55476776665
$ws.Range("A29").Value = 44881.55685109999
$ws.Range("A30").Value = 44881.55893443332
$ws.Range("A31").Value = 44881.56101776665
$ws.Range("A32").Value = 44881.56310109999
$ws.Range("A33").Value = 44881.56518443332
$ws.Range("A34").Value = 44881.56726776665
$ws.Range("A35").Value = 44881.56935109999
$ws.Range("A36").Value = 44881.57143443332
$ws.Range("A37").Value = 44881.57351776666
$ws.Range("A38").Value = 44881.57560109999
$ws.Range("A39").Value = 44881.57768443332
$ws.Range("A40").Value = 44881.57976776666
$ws.Range("A41").Value = 44881.58185109999
$ws.Range("A42").Value = 44881.58393443332
$ws.Range("A43").Value = 44881.58601776665
$ws.Range("A44").Value = 44881.58810109999
$ws.Range("A45").Value = 44881.59018443332
$ws.Range("A46").Value = 44881.59226776665
$ws.Range("A47").Value = 44881.59435109999
$ws.Range("A48").Value = 44881.59643443332
$ws.Range("A49").Value = 44881.59851776665
$ws.Range("A50").Value = 44881.60060109999
$ws.Range("A51").Value = 44881.60268443332
$ws.Range("A52").Value = 44881.60476776666
$ws.Range("A53").Value = 44881.60685109999
$ws.Range("A54").Value = 44881.60893443332
$ws.Range("A55").Value = 44881.61101776666
$ws.Range("A56").Value = 44881.61310109999
$ws.Range("A57").Value = 44881.61518443332
$ws.Range("A58").Value = 44881.61726776665
$ws.Range("A59").Value = 44881.61935109999
$ws.Range("A60").Value = 44881.62143443332
$ws.Range("A61").Value = 44881.62351776665
$ws.Range("A62").Value = 44881.62560109999
$ws.Range("A63").Value = 44881.62768443332
$ws.Range("A64").Value = 44881.62976776665
$ws.Range("A65").Value = 44881.63185109999
$ws.Range("A66").Value = 44881.63393443332
$ws.Range("A67").Value = 44881.63601776666
$ws.Range("A68").Value = 44881.63810109999
$ws.Range("A69").Value = 44881.64018443332
$ws.Range("A70").Value = 44881.64226776666
$ws.Range("A71").Value = 44881.64435109999
$ws.Range("A72").Value = 44881.64643443332
$ws.Range("A73").Value = 44881.64851776665
$ws.Range("A74").Value = 44881.65060109999
$ws.Range("A75").Value = 44881.65268443332
$ws.Range("A76").Value = 44881.65476776665
$ws.Range("A77").Value = 44881.65685109999
$ws.Range("A78").Value = 44881.65893443332
$ws.Range("A79").Value = 44881.66101776665
$ws.Range("A80").Value = 44881.66310109999
$ws.Range("A81").Value = 44881.66518443332
$ws.Range("A82").Value = 44881.66726776666
$ws.Range("A83").Value = 44881.66935109999
$ws.Range("A84").Value = 44881.67143443332
$ws.Range("A85").Value = 44881.67351776666
$ws.Range("A86").Value = 44881.67560109999
$ws.Range("A87").Value = 44881.67768443332
$ws.Range("A88").Value = 44881.67976776665
$ws.Range("A89").Value = 44881.68185109999
$ws.Range("A90").Value = 44881.68393443332
$ws.Range("A91").Value = 44881.68601776665
$ws.Range("A92").Value = 44881.68810109999
$ws.Range("A93").Value = 44881.69018443332
$ws.Range("A94").Value = 44881.69226776665
$ws.Range("A95").Value = 44881.69435109999
$ws.Range("A96").Value = 44881.69643443332
$ws.Range("A97").Value = 44881.69851776666
$ws.Range("A98").Value = 44881.70060109999
$ws.Range("A99").Value = 44881.70268443332
$ws.Range("A100").Value = 44881.70476776666
$ws.Range("A101").Value = 44881.70685109999
$ws.Range("A102").Value = 44881.70893443332
$ws.Range("A103").Value = 44881.71101776665
$ws.Range("A104").Value = 44881.71310109999
$ws.Range("A105").Value = 44881.71518443332
$ws.Range("A106").Value = 44881.71726776665
$ws.Range("A107").Value = 44881.71935109999
$ws.Range("A108").Value = 44881.72143443332
$ws.Range("A109").Value = 44881.72351776665
$ws.Range("A110").Value = 44881.72560109999
$ws.Range("A111").Value = 44881.72768443332
$ws.Range("A112").Value = 44881.72976776666
$ws.Range("A113").Value = 44881.73185109999
$ws.Range("A114").Value = 44881.73393443332
$ws.Range("A115").Value = 44881.73601776666
$ws.Range("A116").Value = 44881.73810109999
$ws.Range("A117").Value = 44881.74018443332
$ws.Range("A118").Value = 44881.74226776665
$ws.Range("A119").Value = 44881.74435109999
$ws.Range("A120").Value = 44881.74643443332
$ws.Range("A121").Value = 44881.74851776665
$ws.Range("A122").Value = 44881.75060109999
